# Major update to combine projects before algorithm
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("P1")
$ws2 = $wb.Worksheets.Item("P2")
$ws3 = $wb.Worksheets.Item("P3")

# Update the Predecessors value for the "FINAL DESIGN REVIEW" row (E6) on each
# of the three project sheets, combining the phases ("3;4;5" -> "2;3;4").
$ws1.Range("E6").Value = "2;3;4"
$ws2.Range("E6").Value = "2;3;4"
$ws3.Range("E6").Value = "2;3;4"

# Update the selection (active cell) remembered on each sheet, and move the
# active/selected tab from P2 to P1.
[void]$ws2.Range("E6").Select()
[void]$ws3.Range("C16").Select()
[void]$ws1.Range("D19").Select()
[void]$ws1.Activate()
